$wb = $excel.ActiveWorkbook

# --- Sheet "Astronauta": add column C values (rows 2-7) ---
$ws = $wb.Worksheets.Item("Astronauta")
$ws.Activate() | Out-Null
$ws.Range("C2").Value = 1
$ws.Range("C3").Value = 1
$ws.Range("C4").Value = 1
$ws.Range("C5").Value = 0.5
$ws.Range("C6").Value = 1
$ws.Range("C7").Value = 1
$ws.Range("C8").Select() | Out-Null

# --- Sheet "Senador": add column C values (rows 2-7) ---
$ws = $wb.Worksheets.Item("Senador")
$ws.Activate() | Out-Null
$ws.Range("C2").Value = 1
$ws.Range("C3").Value = 0.8
$ws.Range("C4").Value = 0.8
$ws.Range("C5").Value = 1
$ws.Range("C6").Value = 1
$ws.Range("C7").Value = 0.8
$ws.Range("C5").Select() | Out-Null

# --- Sheet "Mago": add column C values (rows 2-7) ---
$ws = $wb.Worksheets.Item("Mago")
$ws.Activate() | Out-Null
$ws.Range("C2").Value = 1
$ws.Range("C3").Value = 1
$ws.Range("C4").Value = 1
$ws.Range("C5").Value = 1
$ws.Range("C6").Value = 1
$ws.Range("C7").Value = 1
$ws.Range("C8").Select() | Out-Null

# --- Sheet "Ninja": add column D values (rows 2-7) ---
$ws = $wb.Worksheets.Item("Ninja")
$ws.Activate() | Out-Null
$ws.Range("D2").Value = 1
$ws.Range("D3").Value = 1
$ws.Range("D4").Value = 1
$ws.Range("D5").Value = 1
$ws.Range("D6").Value = 1
$ws.Range("D7").Value = 1
$ws.Range("D8").Select() | Out-Null
